$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 13
$ws.Range("I11").Value = 13
$ws.Range("K11").Value = 13
$ws.Range("M11").Value = 127
$ws.Range("H42").Value = 206
$ws.Range("I42").Value = 235
$ws.Range("J42").Value = 119
$ws.Range("K42").Value = 705
$ws.Range("L42").Value = 357
$ws.Range("M42").Value = -475
$ws.Range("N42").Value = -817
$ws.Range("H130").Value = 65000
$ws.Range("J130").Value = 65000
$ws.Range("L130").Value = 65000
$ws.Range("N130").Value = -75040

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2960.4285
$ws.Range("I45").Value = 2175.8
$ws.Range("J45").Value = 4922
$ws.Range("K45").Value = 2175.8
$ws.Range("L45").Value = 4922
$ws.Range("M45").Value = -1798.8
$ws.Range("N45").Value = -5676
$ws.Range("H61").Value = 2434.7083
$ws.Range("I61").Value = 1207.4375
$ws.Range("K61").Value = 1207.4375
$ws.Range("M61").Value = -995.4375
$ws.Range("H74").Value = 1797.4615
$ws.Range("I74").Value = 1612.8182
$ws.Range("J74").Value = 2813
$ws.Range("K74").Value = 1612.8182
$ws.Range("L74").Value = 2813
$ws.Range("M74").Value = -738.8181999999999
$ws.Range("N74").Value = -4561
$ws.Range("H77").Value = 1797.4615
$ws.Range("I77").Value = 1612.8182
$ws.Range("J77").Value = 2813
$ws.Range("K77").Value = 8064.090999999999
$ws.Range("L77").Value = 14065
$ws.Range("M77").Value = -3696.090999999999
$ws.Range("N77").Value = -22801
$ws.Range("H136").Value = 2434.7083
$ws.Range("I136").Value = 1207.4375
$ws.Range("K136").Value = 3622.3125
$ws.Range("M136").Value = -1072.3125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4505.3335
$ws.Range("I20").Value = 5002
$ws.Range("J20").Value = 4008.6667
$ws.Range("K20").Value = 5002
$ws.Range("L20").Value = 4008.6667
$ws.Range("M20").Value = -4755
$ws.Range("N20").Value = -4502.6667
$ws.Range("H75").Value = 28333.334
$ws.Range("I75").Value = 5000
$ws.Range("J75").Value = 40000
$ws.Range("K75").Value = 5000
$ws.Range("L75").Value = 40000
$ws.Range("M75").Value = -4064
$ws.Range("N75").Value = -41872
$ws.Range("H78").Value = 28333.334
$ws.Range("I78").Value = 5000
$ws.Range("J78").Value = 40000
$ws.Range("K78").Value = 15000
$ws.Range("L78").Value = 120000
$ws.Range("M78").Value = -10320
$ws.Range("N78").Value = -129360
$ws.Range("H86").Value = 5464.6
$ws.Range("J86").Value = 7752.5
$ws.Range("L86").Value = 7752.5
$ws.Range("N86").Value = -9998.5
$ws.Range("H89").Value = 5464.6
$ws.Range("J89").Value = 7752.5
$ws.Range("L89").Value = 38762.5
$ws.Range("N89").Value = -49994.5
$ws.Range("H99").Value = 71429700
$ws.Range("I99").Value = 83334460
$ws.Range("K99").Value = 83334460
$ws.Range("M99").Value = -83332962
$ws.Range("H107").Value = 26320576
$ws.Range("I107").Value = 45456132
$ws.Range("K107").Value = 45456132
$ws.Range("M107").Value = -45454212

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 62.142857
$ws.Range("I7").Value = 41.5
$ws.Range("K7").Value = 41.5
$ws.Range("M7").Value = 71.5
$ws.Range("H99").Value = 4981.8
$ws.Range("J99").Value = 4978.6665
$ws.Range("L99").Value = 4978.6665
$ws.Range("N99").Value = -7974.6665
$ws.Range("H122").Value = 242
$ws.Range("I122").Value = 248.71428
$ws.Range("K122").Value = 746.14284
$ws.Range("M122").Value = 1703.85716
$ws.Range("H126").Value = 4981.8
$ws.Range("J126").Value = 4978.6665
$ws.Range("L126").Value = 14935.9995
$ws.Range("N126").Value = -19875.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 5519.8
$ws.Range("J39").Value = 5519.8
$ws.Range("L39").Value = 16559.4
$ws.Range("N39").Value = -17147.4
$ws.Range("H55").Value = 1992.1765
$ws.Range("I55").Value = 443
$ws.Range("K55").Value = 1329
$ws.Range("M55").Value = -1152
$ws.Range("H60").Value = 1217.3125
$ws.Range("I60").Value = 172.44444
$ws.Range("J60").Value = 2560.7144
$ws.Range("K60").Value = 517.33332
$ws.Range("L60").Value = 7682.1432
$ws.Range("M60").Value = -266.33332
$ws.Range("N60").Value = -8184.1432
$ws.Range("H64").Value = 1053.3334
$ws.Range("J64").Value = 500
$ws.Range("L64").Value = 1500
$ws.Range("N64").Value = -2040
$ws.Range("H67").Value = 1053.3334
$ws.Range("J67").Value = 500
$ws.Range("L67").Value = 1500
$ws.Range("N67").Value = -3372
$ws.Range("H70").Value = 9999.5
$ws.Range("I70").Value = 9999.5
$ws.Range("K70").Value = 29998.5
$ws.Range("M70").Value = -29683.5
$ws.Range("H73").Value = 9999.5
$ws.Range("I73").Value = 9999.5
$ws.Range("K73").Value = 29998.5
$ws.Range("M73").Value = -28906.5
$ws.Range("H75").Value = 8055
$ws.Range("J75").Value = 8055
$ws.Range("L75").Value = 24165
$ws.Range("N75").Value = -26161
$ws.Range("H78").Value = 8055
$ws.Range("J78").Value = 8055
$ws.Range("L78").Value = 72495
$ws.Range("N78").Value = -82479
$ws.Range("H87").Value = 60000
$ws.Range("I87").Value = 60000
$ws.Range("K87").Value = 180000
$ws.Range("M87").Value = -178752
$ws.Range("H90").Value = 60000
$ws.Range("I90").Value = 60000
$ws.Range("K90").Value = 540000
$ws.Range("M90").Value = -533760
$ws.Range("H93").Value = 2225
$ws.Range("I93").Value = 2500
$ws.Range("J93").Value = 2133.3333
$ws.Range("K93").Value = 7500
$ws.Range("L93").Value = 6399.999899999999
$ws.Range("M93").Value = -5628
$ws.Range("N93").Value = -10143.9999
$ws.Range("H107").Value = 644.8889
$ws.Range("I107").Value = 334.8889
$ws.Range("J107").Value = 799.8889
$ws.Range("K107").Value = 1004.6667
$ws.Range("L107").Value = 2399.6667
$ws.Range("M107").Value = 915.3333
$ws.Range("N107").Value = -6239.6667
$ws.Range("H113").Value = 1553.7368
$ws.Range("J113").Value = 1868.5834
$ws.Range("L113").Value = 5605.7502
$ws.Range("N113").Value = -9945.7502
$ws.Range("H114").Value = 2497
$ws.Range("J114").Value = 1994.5
$ws.Range("L114").Value = 5983.5
$ws.Range("N114").Value = -12491.5
$ws.Range("H121").Value = 1304.2354
$ws.Range("I121").Value = 513.8570999999999
$ws.Range("J121").Value = 1857.5
$ws.Range("K121").Value = 1541.5713
$ws.Range("L121").Value = 5572.5
$ws.Range("M121").Value = -231.5712999999998
$ws.Range("N121").Value = -8192.5
$ws.Range("H129").Value = 2314.4546
$ws.Range("J129").Value = 3294.1428
$ws.Range("L129").Value = 9882.428400000001
$ws.Range("N129").Value = -19882.4284
$ws.Range("H131").Value = 2256.7778
$ws.Range("I131").Value = 937.3333
$ws.Range("J131").Value = 2916.5
$ws.Range("K131").Value = 2811.9999
$ws.Range("L131").Value = 8749.5
$ws.Range("M131").Value = 2228.0001
$ws.Range("N131").Value = -18829.5
$ws.Range("H132").Value = 2985.625
$ws.Range("I132").Value = 1975.5
$ws.Range("K132").Value = 17779.5
$ws.Range("M132").Value = -15249.5
$ws.Range("H140").Value = 2250.5557
$ws.Range("I140").Value = 1496.65
$ws.Range("J140").Value = 4404.5713
$ws.Range("K140").Value = 4489.950000000001
$ws.Range("L140").Value = 13213.7139
$ws.Range("M140").Value = 690.0499999999993
$ws.Range("N140").Value = -23573.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2261.2222
$ws.Range("I102").Value = 1471.8334
$ws.Range("J102").Value = 3840
$ws.Range("K102").Value = 1471.8334
$ws.Range("L102").Value = 3840
$ws.Range("M102").Value = 150.1666
$ws.Range("N102").Value = -7084
$ws.Range("H113").Value = 9375
$ws.Range("J113").Value = 9375
$ws.Range("L113").Value = 9375
$ws.Range("N113").Value = -13715
$ws.Range("H126").Value = 3245.5557
$ws.Range("I126").Value = 3018.5
$ws.Range("K126").Value = 9055.5
$ws.Range("M126").Value = -6585.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1819.8
$ws.Range("J22").Value = 1749.75
$ws.Range("L22").Value = 1749.75
$ws.Range("N22").Value = -2339.75
$ws.Range("H27").Value = 1819.8
$ws.Range("J27").Value = 1749.75
$ws.Range("L27").Value = 1749.75
$ws.Range("N27").Value = -1963.75
$ws.Range("H40").Value = 5447.0527
$ws.Range("I40").Value = 5417.353
$ws.Range("K40").Value = 5417.353
$ws.Range("M40").Value = -5281.353
$ws.Range("H55").Value = 1650.6428
$ws.Range("I55").Value = 1394.2858
$ws.Range("K55").Value = 1394.2858
$ws.Range("M55").Value = -1221.2858
$ws.Range("H122").Value = 2990.2666
$ws.Range("I122").Value = 2988.7693
$ws.Range("K122").Value = 8966.3079
$ws.Range("M122").Value = -6516.3079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 30303762
$ws.Range("I107").Value = 30303762
$ws.Range("K107").Value = 90911286
$ws.Range("M107").Value = -90909366
$ws.Range("H122").Value = 5393.125
$ws.Range("I122").Value = 4211.25
$ws.Range("K122").Value = 12633.75
$ws.Range("M122").Value = -10183.75
